# Update gh-pages to output generated at 456a3b4
#
# Changes:
#  1. Sheet "展览": bump several "want to go" counters (col F).
#  2. Sheet "演出": bump F7 counter; insert a new event row (2024-08-02,
#     "广州·井草圣二 2024《夏日独白》指弹吉他音乐会") before the existing
#     2024-08-03 row, shifting the rows below it down by one.
#  3. Sheet "全部类型": same counter bumps + same new event row inserted in
#     the same relative chronological position.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $text) {
    # Plain .Value assignment lets Excel "smart type" a pure ISO date
    # string (e.g. "2024-08-02") into a real date serial. Force it to be
    # treated/stored as text, then drop the leftover explicit number
    # format so the cell ends up with no style override (matching the
    # rest of the sheet, which leaves these cells unstyled).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

function Bump-Counter($ws, $addr, $newValue) {
    $ws.Range($addr).Value = $newValue
}

function Insert-Event($ws, $row, $lastRow, $b, $c, $d, $e, $f, $g, $h, $i) {
    # Shift columns B:I of [row .. lastRow-1] down into [row+1 .. lastRow].
    # Column A (the running index) is intentionally left untouched at each
    # row position -- that matches the source edit, which only shifted the
    # event data and appended one brand-new index cell at the bottom.
    $srcRange = "B" + $row + ":I" + ($lastRow - 1)
    $dstRange = "B" + ($row + 1) + ":I" + $lastRow
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4104)   # xlPasteAll
    $excel.CutCopyMode = $false

    # New trailing index cell, copying A-column's format (bold/border/
    # centered) from the row above it, then writing the correct number.
    $lastA = "A" + ($lastRow - 1)
    $newA = "A" + $lastRow
    $ws.Range($lastA).Copy()
    $ws.Range($newA).PasteSpecial(-4122)       # xlPasteFormats
    $excel.CutCopyMode = $false
    $oldIndexVal = $ws.Range($lastA).Value2
    $ws.Range($newA).Value = $oldIndexVal + 1

    # Fill in the freshly vacated row with the new event's data.
    Set-TextValue $ws ("B" + $row) $b
    $ws.Range("C" + $row).Value = $c
    $ws.Range("D" + $row).Value = $d
    $ws.Range("E" + $row).Value = $e
    $ws.Range("F" + $row).Value = $f
    $ws.Range("G" + $row).Value = $g
    $ws.Range("H" + $row).Value = $h
    $ws.Range("I" + $row).Value = $i
}

# ---------------------------------------------------------------------
# Sheet "展览" -- counter bumps only
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
Bump-Counter $wsExpo "F3"  13856
Bump-Counter $wsExpo "F4"  13637
Bump-Counter $wsExpo "F5"  1054
Bump-Counter $wsExpo "F9"  83
Bump-Counter $wsExpo "F11" 60
Bump-Counter $wsExpo "F13" 2153
Bump-Counter $wsExpo "F22" 329
Bump-Counter $wsExpo "F24" 843

# ---------------------------------------------------------------------
# Sheet "演出" -- counter bump + new row insert
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
Bump-Counter $wsShow "F7" 1543
Insert-Event $wsShow 10 14 `
    "2024-08-02" `
    "广州·井草圣二 2024《夏日独白》指弹吉他音乐会" `
    "恩宁路265号3层 MaoLivehouse(永庆坊店)" `
    "2024.08.02 19:30-08.02 21:00" `
    0 `
    260 `
    "https://show.bilibili.com/platform/detail.html?id=86940" `
    "//i0.hdslb.com/bfs/openplatform/202406/iNGVydXM1717644835981.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" -- counter bumps + new row insert
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Bump-Counter $wsAll "F4"  13856
Bump-Counter $wsAll "F5"  13637
Bump-Counter $wsAll "F6"  1054
Bump-Counter $wsAll "F10" 83
Bump-Counter $wsAll "F12" 60
Bump-Counter $wsAll "F16" 2153
Bump-Counter $wsAll "F29" 329
Bump-Counter $wsAll "F33" 1543
Insert-Event $wsAll 36 42 `
    "2024-08-02" `
    "广州·井草圣二 2024《夏日独白》指弹吉他音乐会" `
    "恩宁路265号3层 MaoLivehouse(永庆坊店)" `
    "2024.08.02 19:30-08.02 21:00" `
    0 `
    260 `
    "https://show.bilibili.com/platform/detail.html?id=86940" `
    "//i0.hdslb.com/bfs/openplatform/202406/iNGVydXM1717644835981.jpeg"
